# Automatische test-sync: 2025-08-03 13:41:50
#
# Appends a new log entry (row 5) to the "Logs" sheet, extends the
# conditional-formatting ranges that cover the data rows so they include
# the new row, and bumps the matching "Aantal" counter on the "Dashboard"
# sheet from 3 to 4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Append the new row to the "Logs" sheet
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 5

$logs.Range("A$newRow").Value = "Kun jij dit even regelen?"
$logs.Range("B$newRow").Value = "mailmind.test@zohomail.eu"
$logs.Range("C$newRow").Value = "Testmail #1: Kun jij dit even regelen?"
$logs.Range("D$newRow").Value = "Intern verzoek / Actie voor medewerker"

$antwoord = @"
Geachte klant,
Dank u wel voor uw bericht. Om u zo goed mogelijk van dienst te kunnen zijn, zou ik graag meer details willen ontvangen over wat u precies geregeld wilt hebben. Kunt u alstublieft specifiëren welke specifieke taak u uitgevoerd wilt hebben?
Met vriendelijke groet,
[Naam] 
E-mailassistent
"@
$logs.Range("E$newRow").Value = $antwoord

$logs.Range("F$newRow").Value = "2025-08-03 13:40:54"
$logs.Range("G$newRow").Value = "Ja"
$logs.Range("H$newRow").Value = "Nee"
$logs.Range("I$newRow").Value = "Ja"
$logs.Range("J$newRow").Value = "Nee"

# The multi-line "Antwoord" text otherwise triggers an automatic custom
# row height; auto-fitting puts the row back to the sheet's default
# height so the row stays free of explicit height attributes, same as
# the other data rows.
$logs.Rows($newRow).AutoFit()

# ---------------------------------------------------------------------
# 2. Extend the conditional formatting sqref ranges (row span 2-4 -> 2-5)
#    for every column that carries conditional formatting.
# ---------------------------------------------------------------------
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "4")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "5")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 3. Update the Dashboard counter for the matching category (3 -> 4)
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 4
